# Auto-generated edit script: updates market-price derived columns (H-N)
# across 8 worksheets per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 81.125
$ws.Range("I8").Value = 64.14286
$ws.Range("K8").Value = 192.42858
$ws.Range("M8").Value = -53.42858000000001

$ws.Range("H32").Value = 920
$ws.Range("I32").Value = 866.6667
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 866.6667
$ws.Range("L32").Value = 1000
$ws.Range("M32").Value = -540.6667
$ws.Range("N32").Value = -1652

$ws.Range("H38").Value = 859.4
$ws.Range("J38").Value = 2000
$ws.Range("L38").Value = 6000
$ws.Range("N38").Value = -6744

$ws.Range("H69").Value = 1548.3334
$ws.Range("J69").Value = 1510.75
$ws.Range("L69").Value = 4532.25
$ws.Range("N69").Value = -6280.25

$ws.Range("H70").Value = 1395.2941
$ws.Range("I70").Value = 1464
$ws.Range("J70").Value = 1297.1428
$ws.Range("K70").Value = 4392
$ws.Range("L70").Value = 3891.4284
$ws.Range("M70").Value = -4122
$ws.Range("N70").Value = -4431.428400000001

$ws.Range("H72").Value = 1548.3334
$ws.Range("J72").Value = 1510.75
$ws.Range("L72").Value = 13596.75
$ws.Range("N72").Value = -22332.75

$ws.Range("H73").Value = 1395.2941
$ws.Range("I73").Value = 1464
$ws.Range("J73").Value = 1297.1428
$ws.Range("K73").Value = 4392
$ws.Range("L73").Value = 3891.4284
$ws.Range("M73").Value = -3456
$ws.Range("N73").Value = -5763.428400000001

$ws.Range("H106").Value = 2937.606
$ws.Range("I106").Value = 1902.0769
$ws.Range("K106").Value = 1902.0769
$ws.Range("M106").Value = -1271.0769

$ws.Range("H118").Value = 437.125
$ws.Range("I118").Value = 437.125
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 1311.375
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 345.625
$ws.Range("N118").ClearContents()

$ws.Range("H121").Value = 6418.75
$ws.Range("J121").Value = 7070.8335
$ws.Range("L121").Value = 21212.5005
$ws.Range("N121").Value = -24706.5005

$ws.Range("H129").Value = 257361.88
$ws.Range("J129").Value = 278759.25
$ws.Range("L129").Value = 836277.75
$ws.Range("N129").Value = -846277.75

$ws.Range("H132").Value = 20410270
$ws.Range("I132").Value = 23258040
$ws.Range("J132").Value = 1252.6666
$ws.Range("K132").Value = 69774120
$ws.Range("L132").Value = 3757.9998
$ws.Range("M132").Value = -69771590
$ws.Range("N132").Value = -8817.9998

$ws.Range("H137").Value = 107933.69
$ws.Range("I137").Value = 127118.125
$ws.Range("K137").Value = 381354.375
$ws.Range("M137").Value = -378804.375

$ws.Range("H138").Value = 2899.863
$ws.Range("J138").Value = 3154.3333
$ws.Range("L138").Value = 9462.999899999999
$ws.Range("N138").Value = -19742.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13863.702
$ws.Range("I32").Value = 10171.5
$ws.Range("K32").Value = 10171.5
$ws.Range("M32").Value = -9884.5

$ws.Range("H88").Value = 144451.58
$ws.Range("I88").Value = 1850
$ws.Range("J88").Value = 201492.2
$ws.Range("K88").Value = 1850
$ws.Range("L88").Value = 201492.2
$ws.Range("M88").Value = -1444
$ws.Range("N88").Value = -202304.2

$ws.Range("H91").Value = 144451.58
$ws.Range("I91").Value = 1850
$ws.Range("J91").Value = 201492.2
$ws.Range("K91").Value = 1850
$ws.Range("L91").Value = 201492.2
$ws.Range("M91").Value = -446
$ws.Range("N91").Value = -204300.2

$ws.Range("H102").Value = 1463.6316
$ws.Range("I102").Value = 1461.6111
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 1461.6111
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 160.3888999999999
$ws.Range("N102").Value = -4744

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1477.7778
$ws.Range("J99").Value = 1050
$ws.Range("L99").Value = 1050
$ws.Range("N99").Value = -4046

$ws.Range("H134").Value = 4462.6387
$ws.Range("I134").Value = 4295.25
$ws.Range("J134").Value = 5048.5
$ws.Range("K134").Value = 12885.75
$ws.Range("L134").Value = 15145.5
$ws.Range("M134").Value = -10350.75
$ws.Range("N134").Value = -20215.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 175
$ws.Range("I22").Value = 175
$ws.Range("K22").Value = 175
$ws.Range("M22").Value = 175

$ws.Range("H31").Value = 5137.93
$ws.Range("I31").Value = 2560.8823
$ws.Range("K31").Value = 2560.8823
$ws.Range("M31").Value = -2265.8823

$ws.Range("H34").Value = 5137.93
$ws.Range("I34").Value = 2560.8823
$ws.Range("K34").Value = 2560.8823
$ws.Range("M34").Value = -2358.8823

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 106
$ws.Range("I6").Value = 74.833336
$ws.Range("K6").Value = 224.500008
$ws.Range("M6").Value = -111.500008

$ws.Range("H7").Value = 1675.5
$ws.Range("I7").Value = 1675.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 5026.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4914.5
$ws.Range("N7").ClearContents()

$ws.Range("H23").Value = 276.17392
$ws.Range("J23").Value = 315.1
$ws.Range("L23").Value = 945.3000000000001
$ws.Range("N23").Value = -1415.3

$ws.Range("H107").Value = 6717.533
$ws.Range("I107").Value = 8927.817999999999
$ws.Range("J107").Value = 639.25
$ws.Range("K107").Value = 26783.454
$ws.Range("L107").Value = 1917.75
$ws.Range("M107").Value = -24863.454
$ws.Range("N107").Value = -5757.75

$ws.Range("H122").Value = 1358.7241
$ws.Range("I122").Value = 590
$ws.Range("J122").Value = 1415.6666
$ws.Range("K122").Value = 5310
$ws.Range("L122").Value = 12740.9994
$ws.Range("M122").Value = -2860
$ws.Range("N122").Value = -17640.9994

$ws.Range("H131").Value = 722.9
$ws.Range("I131").Value = 405.83334
$ws.Range("J131").Value = 766.13635
$ws.Range("K131").Value = 1217.50002
$ws.Range("L131").Value = 2298.40905
$ws.Range("M131").Value = 3822.49998
$ws.Range("N131").Value = -12378.40905

$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4962.5
$ws.Range("I102").Value = 4801.375
$ws.Range("K102").Value = 4801.375
$ws.Range("M102").Value = -3179.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1511203.5
$ws.Range("I122").Value = 2453205
$ws.Range("J122").Value = 4001
$ws.Range("K122").Value = 7359615
$ws.Range("L122").Value = 12003
$ws.Range("M122").Value = -7357165
$ws.Range("N122").Value = -16903

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 35000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 35000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 35000
$ws.Range("N54").Value = -36040
$ws.Range("M54").ClearContents()
